# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K") held strikeout counts that were previously computed
# against an older "Strike#" definition. The stat pipeline was regenerated to
# source K from the updated field, and the recalculated std/mean derived
# s_vals are written back into column G for every data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values recalculated by the regenerated save_data /
# std-mean pipeline, keyed by worksheet row number.
$sVals = [ordered]@{
    2  = 0
    3  = 1
    4  = 1
    6  = 0
    7  = 0
    8  = 2
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 2
    16 = 0
    17 = 2
    18 = 0
    19 = 1
    20 = 2
    21 = 1
    22 = 1
    23 = 2
    24 = 2
    25 = 2
    26 = 0
    27 = 1
    28 = 1
    30 = 1
    31 = 2
    32 = 0
    33 = 1
    34 = 0
    35 = 1
    36 = 1
    38 = 2
    39 = 0
    40 = 0
    41 = 0
    42 = 1
    43 = 1
    44 = 1
    45 = 2
    47 = 2
    48 = 1
    49 = 1
}

foreach ($row in $sVals.Keys) {
    $ws.Range("G$row").Value = $sVals[$row]
}
